# Add the "11thSep, 2021" daily log sheet (new last tab, sheetId 9) with the
# same layout used by every other daily sheet in this workbook:
#   A1 = "Name"
#   A2 = "BARKAT", B2 = "Time",        C2:F2 = timestamp serials
#                  B3 = "Matching(%)", C3:F3 = matching percentages

$wb = $excel.ActiveWorkbook

# The workbook's style table already carries two time-only formats
# (numFmtId 164/165). Register the new "yyyy-mm-dd h:mm:ss" number format
# used elsewhere in this tracker via a scratch sheet so it lands in
# styles.xml without changing the formatting actually applied to any
# worksheet cell, then discard the scratch sheet.
$scratch = $wb.Worksheets.Add()
$scratch.Range("A1").NumberFormat = "yyyy-mm-dd h:mm:ss"
[void]$scratch.Delete()

# Duplicate the most recent daily sheet ("03thSep, 2021") so the new tab
# inherits its exact layout/styles (bold labels, time-formatted cells,
# page margins, etc.), then drop it in right after as the new last tab.
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$lastSheet.Copy($null, $lastSheet)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "11thSep, 2021"

# Overwrite the copied Time / Matching(%) rows with this day's figures.
$ws.Range("C2").Value = 44450.87703130787
$ws.Range("D2").Value = 44450.87708763889
$ws.Range("E2").Value = 44450.87713575231
$ws.Range("F2").Value = 44450.87718730697

$ws.Range("C3").Value = 57.38756673793869
$ws.Range("D3").Value = 56.41789086541215
$ws.Range("E3").Value = 54.25104163473586
$ws.Range("F3").Value = 51.08050264201491
